$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the new daily-report rows (24-30) for 2025-12-20 .. 2025-12-26.
#    Formats are cloned from existing rows with the matching style pattern
#    (via PasteSpecial xlPasteFormats = -4122) so no new style entries are
#    created in styles.xml, then the cell values are written on top.
# ---------------------------------------------------------------------------

# Row 24 (Sat 46011) and Row 25 (Sun 46012): same pattern as rows 16 / 17
$ws.Range("A16:E16").Copy() | Out-Null
$ws.Range("A24:E24").PasteSpecial(-4122) | Out-Null
$ws.Range("A25:E25").PasteSpecial(-4122) | Out-Null

$ws.Range("A24").Value = 46011
$ws.Range("B24").Value = "-"
$ws.Range("C24").Value = "holiday"
$ws.Range("D24").Value = "holiday"
$ws.Range("E24").Value = "-"

$ws.Range("A25").Value = 46012
$ws.Range("B25").Value = "-"
$ws.Range("C25").Value = "holiday"
$ws.Range("D25").Value = "holiday"
$ws.Range("E25").Value = "-"

# Row 26 (Mon 46013): same pattern as row 23 (A:F)
$ws.Range("A23:F23").Copy() | Out-Null
$ws.Range("A26:F26").PasteSpecial(-4122) | Out-Null

$ws.Range("A26").Value = 46013
$ws.Range("B26").Value = 16
$ws.Range("C26").Value = "Express and JSON (we need to store a data on json file on public folderder)"
$ws.Range("D26").Value = "Express and json"
$ws.Range("F26").Value = "D:\intership\task\third_week\task_3_Blog_on_react\Intership_task\third_week\task_3_Blog_on_react\blog-app"

# Row 27 (Tue 46014): same pattern as row 14 (A:D)
$ws.Range("A14:D14").Copy() | Out-Null
$ws.Range("A27:D27").PasteSpecial(-4122) | Out-Null

$ws.Range("A27").Value = 46014
$ws.Range("B27").Value = 17
$ws.Range("C27").Value = "Express and JSON (we need to store a data on json file on public folderder)"
$ws.Range("D27").Value = "Express and json"

# Rows 28-30: date + day number only, same pattern as row 4 (A:B)
$ws.Range("A4:B4").Copy() | Out-Null
$ws.Range("A28:B28").PasteSpecial(-4122) | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:B30").PasteSpecial(-4122) | Out-Null

$ws.Range("A28").Value = 46015
$ws.Range("B28").Value = 18

$ws.Range("A29").Value = 46016
$ws.Range("B29").Value = 19

$ws.Range("A30").Value = 46017
$ws.Range("B30").Value = 20

# Row heights to match the rest of the table (15.6 points).
$ws.Range("A24:H30").RowHeight = 15.6

# Hyperlink for the new task entry on row 26.
$ws.Hyperlinks.Add($ws.Range("E26"), "https://github.com/AakashChidambaranathan/Intership_task/tree/2385d676f720ebd8a58de01590373056baf3a957", "", "", "AakashChidambaranathan/Intership_task at 2385d676f720ebd8a58de01590373056baf3a957") | Out-Null

# ---------------------------------------------------------------------------
# 2. Swap the H-column "open item" placeholder styles (cosmetic re-order of
#    cellXfs entries 19 / 20 that Excel performed when it re-saved the file,
#    same rendered look, different xf index). Captured via scratch cells so
#    the paste-special calls reuse the existing style slots instead of
#    creating new ones.
# ---------------------------------------------------------------------------
$ws.Range("H4").Copy() | Out-Null
$ws.Range("AZ100").PasteSpecial(-4122) | Out-Null

$ws.Range("H12").Copy() | Out-Null
$ws.Range("AZ101").PasteSpecial(-4122) | Out-Null

$ws.Range("AZ101").Copy() | Out-Null
$ws.Range("H4:H7").PasteSpecial(-4122) | Out-Null

$ws.Range("AZ100").Copy() | Out-Null
$ws.Range("H12:H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H20:H23").PasteSpecial(-4122) | Out-Null

$ws.Range("AZ100").Clear() | Out-Null
$ws.Range("AZ101").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3. Window/view adjustments: scroll to A2 and change the zoom level,
#    move the active selection to C33.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 93
$ws.Range("C33").Select() | Out-Null
